$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is numeric-looking must be forced to Text so they
# keep matching the source data (all Price/Volume cells are stored as text).
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

$ws.Range('D2').Value = '58.377.22'
$ws.Range('E2').Value = '  -2.94%  '
$ws.Range('D3').Value = '2.450.91'
$ws.Range('E3').Value = '  -3.87%  '
$ws.Range('E4').Value = '  +0.06%  '
Set-TextValue 'D5' '528.32'
$ws.Range('E5').Value = '  -2.22%  '
Set-TextValue 'D6' '133.85'
$ws.Range('E6').Value = '  -7.14%  '
Set-TextValue 'D7' '0.999'
$ws.Range('E7').Value = '  +0.22%  '
$ws.Range('E8').Value = '  -2.63%  '
$ws.Range('D9').Value = '2.458.17'
$ws.Range('E9').Value = '  -4.33%  '
Set-TextValue 'D10' '0.0991'
$ws.Range('E10').Value = '  -2.07%  '
$ws.Range('E11').Value = '  -0.35%  '
Set-TextValue 'D12' '5.30'
$ws.Range('E12').Value = '  -3.28%  '
Set-TextValue 'D13' '0.342'
$ws.Range('E13').Value = '  -5.64%  '
$ws.Range('D14').Value = '2.890.72'
$ws.Range('E14').Value = '  -3.57%  '
$ws.Range('D15').Value = '58.296.02'
$ws.Range('E15').Value = '  -2.99%  '
Set-TextValue 'D16' '22.61'
$ws.Range('E16').Value = '  -6.06%  '
Set-TextValue 'D17' '0.0000137'
$ws.Range('E17').Value = '  -4.10%  '
$ws.Range('D18').Value = '2.460.99'
$ws.Range('E18').Value = '  -3.41%  '
Set-TextValue 'D19' '10.72'
$ws.Range('E19').Value = '  -4.68%  '
$ws.Range('E20').Value = '  -3.01%  '
$ws.Range('E21').Value = '  -1.99%  '
$ws.Range('E22').Value = '  -0.24%  '
$ws.Range('E23').Value = '  -4.04%  '
Set-TextValue 'D24' '62.43'
$ws.Range('E24').Value = '  -1.74%  '
$ws.Range('E25').Value = '  -5.99%  '
$ws.Range('E26').Value = '  -1.95%  '
$ws.Range('E27').Value = '  -1.07%  '
$ws.Range('E28').Value = '  -7.65%  '
$ws.Range('D29').Value = '0.0₃0751'
$ws.Range('E29').Value = '  -5.52%  '
Set-TextValue 'D30' '6.50'
$ws.Range('E30').Value = '  -7.98%  '
$ws.Range('E31').Value = '  -3.69%  '
Set-TextValue 'D32' '164.11'
$ws.Range('E32').Value = '  -0.88%  '
Set-TextValue 'D33' '0.998'
$ws.Range('E33').Value = '  +0.04%  '
$ws.Range('E34').Value = '  -7.20%  '
$ws.Range('E35').Value = '  -8.45%  '
$ws.Range('E36').Value = '  -2.85%  '
Set-TextValue 'D37' '4.01'
$ws.Range('E37').Value = '  -8.68%  '
Set-TextValue 'D38' '1.53'
$ws.Range('E38').Value = '  -5.99%  '
Set-TextValue 'D39' '36.46'
$ws.Range('E39').Value = '  -1.57%  '
Set-TextValue 'D40' '0.804'
$ws.Range('E40').Value = '  -3.53%  '
$ws.Range('E42').Value = '  -9.34%  '
Set-TextValue 'D43' '273.90'
$ws.Range('E43').Value = '  -9.25%  '
Set-TextValue 'D44' '0.998'
$ws.Range('E44').Value = '  +0.36%  '
$ws.Range('E45').Value = '  -0.29%  '
Set-TextValue 'D46' '0.585'
$ws.Range('E46').Value = '  -4.31%  '
$ws.Range('E47').Value = '  -1.66%  '
Set-TextValue 'D48' '120.76'
$ws.Range('E48').Value = '  -5.09%  '
Set-TextValue 'D49' '0.0504'
$ws.Range('E49').Value = '  -2.79%  '
Set-TextValue 'D50' '0.0217'
$ws.Range('E50').Value = '  -4.97%  '
Set-TextValue 'D51' '17.07'
$ws.Range('E51').Value = '  -6.44%  '
